{"js": "// Replace the twenty-five three-digit-by-one-digit multiplication prompts\n// with their new values. Each old expression is unique in the document, so\n// a plain text search/replace (matchCase, whole-document body search) is\n// unambiguous. The pairs are applied in the same order they appear in the\n// source document so that a value that is simultaneously an old prompt\n// (elsewhere) and a new prompt (here) \u2014 \"602\u00d76=\" \u2014 is never re-matched\n// after it has been freshly inserted.\nconst replacements = [\n  [\"703\u00d73=\", \"500\u00d75=\"],\n  [\"222\u00d73=\", \"586\u00d72=\"],\n  [\"800\u00d72=\", \"341\u00d74=\"],\n  [\"144\u00d77=\", \"945\u00d76=\"],\n  [\"152\u00d79=\", \"430\u00d78=\"],\n  [\"840\u00d74=\", \"916\u00d78=\"],\n  [\"430\u00d72=\", \"808\u00d79=\"],\n  [\"644\u00d74=\", \"463\u00d79=\"],\n  [\"183\u00d73=\", \"872\u00d74=\"],\n  [\"636\u00d74=\", \"730\u00d73=\"],\n  [\"239\u00d75=\", \"812\u00d73=\"],\n  [\"295\u00d75=\", \"746\u00d73=\"],\n  [\"677\u00d79=\", \"552\u00d75=\"],\n  [\"821\u00d76=\", \"553\u00d73=\"],\n  [\"759\u00d79=\", \"886\u00d76=\"],\n  [\"365\u00d78=\", \"570\u00d79=\"],\n  [\"772\u00d78=\", \"691\u00d73=\"],\n  [\"602\u00d76=\", \"406\u00d78=\"],\n  [\"742\u00d72=\", \"151\u00d73=\"],\n  [\"270\u00d73=\", \"163\u00d76=\"],\n  [\"891\u00d78=\", \"408\u00d76=\"],\n  [\"886\u00d75=\", \"722\u00d79=\"],\n  [\"156\u00d72=\", \"230\u00d75=\"],\n  [\"522\u00d79=\", \"246\u00d74=\"],\n  [\"149\u00d74=\", \"602\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Expected to find \"${oldText}\" in the document, but it was not present.`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the twenty-five three-digit-by-one-digit multiplication prompts\n# with their new values. Each old expression is unique in the document, so\n# Find/Replace against the whole-document range is unambiguous. Pairs are\n# applied in the same order they appear in the source document so that a\n# value that is simultaneously an old prompt (elsewhere) and a new prompt\n# (here) -- \"602x6=\" -- is never re-matched after it has been freshly\n# inserted by an earlier step.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"703\u00d73=\", \"500\u00d75=\"),\n  @(\"222\u00d73=\", \"586\u00d72=\"),\n  @(\"800\u00d72=\", \"341\u00d74=\"),\n  @(\"144\u00d77=\", \"945\u00d76=\"),\n  @(\"152\u00d79=\", \"430\u00d78=\"),\n  @(\"840\u00d74=\", \"916\u00d78=\"),\n  @(\"430\u00d72=\", \"808\u00d79=\"),\n  @(\"644\u00d74=\", \"463\u00d79=\"),\n  @(\"183\u00d73=\", \"872\u00d74=\"),\n  @(\"636\u00d74=\", \"730\u00d73=\"),\n  @(\"239\u00d75=\", \"812\u00d73=\"),\n  @(\"295\u00d75=\", \"746\u00d73=\"),\n  @(\"677\u00d79=\", \"552\u00d75=\"),\n  @(\"821\u00d76=\", \"553\u00d73=\"),\n  @(\"759\u00d79=\", \"886\u00d76=\"),\n  @(\"365\u00d78=\", \"570\u00d79=\"),\n  @(\"772\u00d78=\", \"691\u00d73=\"),\n  @(\"602\u00d76=\", \"406\u00d78=\"),\n  @(\"742\u00d72=\", \"151\u00d73=\"),\n  @(\"270\u00d73=\", \"163\u00d76=\"),\n  @(\"891\u00d78=\", \"408\u00d76=\"),\n  @(\"886\u00d75=\", \"722\u00d79=\"),\n  @(\"156\u00d72=\", \"230\u00d75=\"),\n  @(\"522\u00d79=\", \"246\u00d74=\"),\n  @(\"149\u00d74=\", \"602\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $ok = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $ok) {\n        throw \"Expected to find '$oldText' in the document, but it was not present.\"\n    }\n}\n"}
